$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.833.09'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.647.71'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").Value = '  +0.44%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.24'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.503'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.252'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0629'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.23'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0843'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.871.57'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.633.86'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.18'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.530'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.79'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.824.51'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0740'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.59'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.50%  '
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.41'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +13.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.31'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.39'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.15'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.33%  '
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("E27").Value = '  -1.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.12'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.74'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.08%  '
$ws.Range("E30").Value = '  -1.04%  '
$ws.Range("E31").Value = '  +0.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.33'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.01'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.297.32'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.66%  '
$ws.Range("E35").Value = '  -0.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.44'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0175'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.541'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.831'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.24%  '
$ws.Range("E40").Value = '  +0.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.811'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.24'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.37'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.795.35'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.71'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.43'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.61'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0519'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0100'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -5.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.71'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("E51").Value = '  +0.09%  '
